# issue #5: stock data from json to db
#
# Adds three new columns (category, source_file, index) to the "股票"
# (stock) sheet (sheet 3) and shifts the existing date / legislator_name /
# legislator_id columns one place to the right to make room for the new
# "category" column right after "property_category".
#
# Resulting column layout (row 1 header / rows 2-14 data):
#   A index(orig)  B name            C owner            D quantity
#   E face_value   F currency        G total            H property_category
#   I category(NEW)  J date          K legislator_name  L legislator_id(NEW)
#   M source_file(NEW)               N index(NEW)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$firstDataRow = 2
$lastDataRow = 14

# ------------------------------------------------------------------
# Header row (row 1)
# ------------------------------------------------------------------

# Give the 3 brand-new header cells (L1:N1) the same look as the existing
# header cells (copy style+border+font from K1).
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

# Shift the existing header text one column to the right using a
# values-only paste (so Excel does not try to reinterpret the text, e.g.
# turn "2011-11-17" into a date serial number). Process right-to-left so
# we never overwrite a source cell before it has been copied.
$ws.Cells.Item(1, 11).Copy()                      # old K1 = legislator_id
$ws.Cells.Item(1, 12).PasteSpecial(-4163)          # -> L1

$ws.Cells.Item(1, 10).Copy()                       # old J1 = legislator_name
$ws.Cells.Item(1, 11).PasteSpecial(-4163)          # -> K1

$ws.Cells.Item(1, 9).Copy()                        # old I1 = date
$ws.Cells.Item(1, 10).PasteSpecial(-4163)          # -> J1

# Fill the new header labels (in the same relative order the columns were
# introduced: category, then source_file, then index).
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# ------------------------------------------------------------------
# Data rows (rows 2-14)
# ------------------------------------------------------------------
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {

    # Remember the row's original "index" value (column A) before we touch
    # anything else.
    $indexVal = $ws.Cells.Item($r, 1).Value()

    # New numeric cells (L = legislator_id, N = index) need the same style
    # as the existing numeric data cells (col A/D/E/G all share style 2... 
    # except col A which is style 1, so copy from col D which is a plain
    # numeric data cell).
    $ws.Cells.Item($r, 4).Copy()
    $ws.Cells.Item($r, 12).PasteSpecial(-4122)
    $ws.Cells.Item($r, 14).PasteSpecial(-4122)

    # New text cell (M = source_file) needs the same style as the existing
    # text data cells (col H).
    $ws.Cells.Item($r, 8).Copy()
    $ws.Cells.Item($r, 13).PasteSpecial(-4122)

    # Shift the existing date / legislator_name / legislator_id values one
    # column to the right with a values-only paste (keeps "2011-11-17" as
    # plain text instead of letting it get parsed into a date). Go
    # right-to-left so sources are read before being overwritten.
    $ws.Cells.Item($r, 11).Copy()                  # old K = legislator_id
    $ws.Cells.Item($r, 12).PasteSpecial(-4163)     # -> L

    $ws.Cells.Item($r, 10).Copy()                  # old J = legislator_name
    $ws.Cells.Item($r, 11).PasteSpecial(-4163)     # -> K

    $ws.Cells.Item($r, 9).Copy()                   # old I = date
    $ws.Cells.Item($r, 10).PasteSpecial(-4163)     # -> J

    # Fill the new columns (again mirroring the introduction order used on
    # the header: category, source_file, index).
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmpbbad1"
    $ws.Cells.Item($r, 14).Value = $indexVal
}

$excel.CutCopyMode = 0
